# Correct the mismatched "incorrect lists" entries on Tabelle1.
# (commit message: "corrected incorrect lists ;)")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

$ws.Range("C1").Value = "se"
$ws.Range("D1").Value = "do"

$ws.Range("B4").Value = "ka"
$ws.Range("C4").Value = "fi"
$ws.Range("G4").Value = "ba"

$ws.Range("D5").Value = "ki"
$ws.Range("E5").Value = "fi"

$ws.Range("C6").Value = "lo"

$ws.Range("C11").Value = "fi"
$ws.Range("F11").Value = "lo"

$ws.Range("D14").Value = "fu"
$ws.Range("E14").Value = "to"
$ws.Range("F14").Value = "ba"
$ws.Range("G14").Value = "fu"

$ws.Range("C18").Value = "pe"
$ws.Range("E18").Value = "lo"
$ws.Range("F18").Value = "ba"
$ws.Range("G18").Value = "ba"

$ws.Range("B25").Value = "mi"
$ws.Range("C25").Value = "la"
$ws.Range("E25").Value = "se"
$ws.Range("F25").Value = "fi"

$ws.Range("E26").Value = "fi"
$ws.Range("F26").Value = "se"

$ws.Range("B27").Value = "lo"
$ws.Range("E27").Value = "se"

$ws.Range("E29").Value = "se"
$ws.Range("F29").Value = "fi"
$ws.Range("G29").Value = "ba"

$ws.Range("B32").Value = "fi"
$ws.Range("C32").Value = "fu"
$ws.Range("D32").Value = "lo"
$ws.Range("F32").Value = "lo"
$ws.Range("G32").Value = "fi"

$ws.Range("D33").Value = "fu"
$ws.Range("F33").Value = "bo"

$ws.Range("B34").Value = "pe"
$ws.Range("C34").Value = "li"
$ws.Range("E34").Value = "ba"
$ws.Range("F34").Value = "pe"

$ws.Range("C35").Value = "lo"
$ws.Range("E35").Value = "pe"
$ws.Range("F35").Value = "fi"

$ws.Range("B40").Value = "fu"
$ws.Range("E40").Value = "pe"

# Reflect the author's final selection/scroll position on the sheet.
$ws.Range("C30").Select()
